$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the name in A3 ("Adrian Geanta" -> "Popescu Ion")
$ws.Range("A3").Value = "Popescu Ion"

# Clear the names that used to be in A4, A5, A6 ("Cosmin Geanta", "Matei Marius", "Matei Ioana")
$ws.Range("A4").Value = ""
$ws.Range("A5").Value = ""
$ws.Range("A6").Value = ""

# Move the active selection from A11 to A5
$ws.Range("A5").Select()
